$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.38999999999999
$ws.Range("B4").Value = 4.762900000000005
$ws.Range("C6").Value = -10.9932
$ws.Range("B7").Value = 6.054700000000004
$ws.Range("C7").Value = -11.0395
$ws.Range("B8").Value = 5.475999999999995
$ws.Range("C8").Value = -10.75839999999999
$ws.Range("A11").Value = -22.18700000000002
$ws.Range("E11").Value = 12.7828
$ws.Range("A12").Value = -22.65460000000001
$ws.Range("B12").Value = 6.091700000000004
$ws.Range("B14").Value = 9.462500000000007
$ws.Range("E14").Value = 12.9537
$ws.Range("A15").Value = -21.30810000000002
$ws.Range("C19").Value = -13.14509999999999
$ws.Range("E19").Value = 12.85529999999999
$ws.Range("C21").Value = -13.1744
$ws.Range("E21").Value = 12.66039999999999
$ws.Range("B22").Value = 5.048200000000005
$ws.Range("C24").Value = -11.2975
$ws.Range("C25").Value = -11.12989999999999
